# Apply the "Changed some datasets and added new policy file" edit:
# append two new country rows (Bulgaria, Slovakia) to the "Data" sheet,
# and move the active selection to the newly added area.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Activate()

# Row 29 - Bulgaria
$ws.Range("A29").Value = "Bulgaria"
$ws.Range("B29").Value = 73.650999999999996
$ws.Range("C29").Value = 73.989999999999995
$ws.Range("D29").Value = 74.328999999999994
$ws.Range("E29").Value = 74.668999999999997
$ws.Range("F29").Value = 75.007999999999996
$ws.Range("G29").Value = 75.346999999999994
$ws.Range("H29").Value = 75.686000000000007
$ws.Range("I29").Value = 76.025000000000006
$ws.Range("J29").Value = 76.363

# Row 30 - Slovakia
$ws.Range("A30").Value = "Slovakia"
$ws.Range("B30").Value = 54
$ws.Range("C30").Value = 53.889000000000003
$ws.Range("D30").Value = 53.805999999999997
$ws.Range("E30").Value = 53.750999999999998
$ws.Range("F30").Value = 53.725999999999999
$ws.Range("G30").Value = 53.728999999999999
$ws.Range("H30").Value = 53.76
$ws.Range("I30").Value = 53.82
$ws.Range("J30").Value = 53.908999999999999

# Update the view selection to match the edited area (mirrors the
# workbook having been scrolled down to the newly-added rows and D29
# being the last-active cell before save).
$ws.Range("D29").Select()

Write-Host "Applied urbanization dataset update (Bulgaria, Slovakia rows added)."
